$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 7
$ws.Range("C8").Value = 20.441999435424805
$ws.Range("D8").Value = 14
